# MALM architecture.pptx — remove the encoder/decoder detour.
#
# The "responselist -> announce -> synthesize" side path duplicated what the
# "decoder/responselist/encoder" column already shows, so it is dropped, the
# two rectangles that fed into it are renamed to what they actually are, and
# the big outer container shrinks to fit the now-narrower diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1) Shrink the big background rectangle ("직사각형 3") now that the
#    right-hand announce/synthesize branch is gone. Only the width changes
#    (7731377 EMU -> 5902577 EMU); Height/Left/Top stay put.
#    Shape.Width is a COM Single (float32 points), so feed it a value whose
#    float32 round-trip lands exactly on the target EMU instead of the
#    naive emu/12700 (which truncates 1 EMU low).
$container = $s.Shapes.Item("직사각형 3")
$container.Width = 464.7698821997165

# 2) Relabel the pipeline boxes: the decode/encode naming was wrong — this
#    is really just building/sending the response list.
$s.Shapes.Item("직사각형 5").TextFrame.TextRange.Text = "responselist"
$s.Shapes.Item("직사각형 23").TextFrame.TextRange.Text = "announce"
$s.Shapes.Item("직사각형 24").TextFrame.TextRange.Text = "synthesize"

# 3) Delete the now-redundant "announce"/"synthesize" boxes and their two
#    connector arrows out past the right edge of the diagram.
$s.Shapes.Item("직사각형 41").Delete()
$s.Shapes.Item("직사각형 42").Delete()
$s.Shapes.Item("직선 화살표 연결선 43").Delete()
$s.Shapes.Item("직선 화살표 연결선 44").Delete()
